$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "CasesTab" query (B2): the `Cohort` clause / column was removed from
# the RETURN statement (trailing comma on the Response to Treatment line was
# also dropped since it is now the last item).
$casesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Belgian Malinois', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in [ 'Bladder, Prostate', 'Bladder, Urethra', 'Bladder, Urethra, Prostate', 'Urethra, Prostate'] and diag.best_response in ['Not Determined']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Cells.Item(2, 2).Value = $casesTabQuery

# Row 2 shrinks now that the query text is shorter.
$ws.Rows.Item(2).RowHeight = 319

# Selection/scroll position moves back up to the top of the sheet.
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Range("B2").Select()
